$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 (Caso 3154), which shifts all subsequent rows up by one.
$ws.Rows.Item(4).Delete()
